# injection results - added average air pressure
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-5 had picked up a stray "apply fill" flag on their style (A3:A5 = style
# w/ applyFill, B3:B5 = default style w/ applyFill). Re-apply the clean formats
# used by the neighboring header cells (A1 for labels in column A, B2 for the
# un-styled values in column B) so those extra cell-format records go away.
$ws.Range("A1").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B2").Copy()
$ws.Range("B3:B5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Insert a new summary row just above "Ave Temp c (water)" (row 34) for the
# newly-tracked average air pressure, shifting everything below down by one.
$ws.Rows.Item(34).Insert()
$ws.Range("A34").Value = "Ave Pressure kpa (air)"
$ws.Range("B34").Value = 63.28

# Leave the selection where the editor ended up.
$ws.Range("B33").Select()
